# Update the Pokemon OU usage stats table: refresh the scraped names/percentages
# for rows 2-69, trim the trailing rows that no longer exist (was B73, now B69),
# and widen column A to match column B's width.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Landorus-Therian"
$ws.Range("B2").Value = 0.2997
$ws.Range("A3").Value = "Garchomp"
$ws.Range("B3").Value = 0.2992
$ws.Range("A4").Value = "Ferrothorn"
$ws.Range("B4").Value = 0.2758
$ws.Range("A5").Value = "Zapdos"
$ws.Range("B5").Value = 0.2692
$ws.Range("A6").Value = "Dragapult"
$ws.Range("B6").Value = 0.2591
$ws.Range("A7").Value = "Weavile"
$ws.Range("B7").Value = 0.2582
$ws.Range("A8").Value = "Heatran"
$ws.Range("B8").Value = 0.2467
$ws.Range("A9").Value = "Clefable"
$ws.Range("B9").Value = 0.2265
$ws.Range("A10").Value = "Rillaboom"
$ws.Range("B10").Value = 0.1974
$ws.Range("A11").Value = "Tapu Lele"
$ws.Range("B11").Value = 0.1904
$ws.Range("A12").Value = "Tornadus-Therian"
$ws.Range("B12").Value = 0.1869
$ws.Range("A13").Value = "Urshifu-Rapid-Strike"
$ws.Range("B13").Value = 0.1831
$ws.Range("A14").Value = "Melmetal"
$ws.Range("B14").Value = 0.1767
$ws.Range("A15").Value = "Kartana"
$ws.Range("B15").Value = 0.1646
$ws.Range("A16").Value = "Slowbro"
$ws.Range("B16").Value = 0.1569
$ws.Range("A17").Value = "Tyranitar"
$ws.Range("B17").Value = 0.1278
$ws.Range("A18").Value = "Excadrill"
$ws.Range("B18").Value = 0.117
$ws.Range("A19").Value = "Toxapex"
$ws.Range("B19").Value = 0.1166
$ws.Range("A20").Value = "Corviknight"
$ws.Range("B20").Value = 0.1151
$ws.Range("A21").Value = "Hatterene"
$ws.Range("B21").Value = 0.1102
$ws.Range("A22").Value = "Blaziken"
$ws.Range("B22").Value = 0.1064
$ws.Range("A23").Value = "Zeraora"
$ws.Range("B23").Value = 0.08529999999999999
$ws.Range("A24").Value = "Slowking-Galar"
$ws.Range("B24").Value = 0.08
$ws.Range("A25").Value = "Tapu Koko"
$ws.Range("B25").Value = 0.0784
$ws.Range("A26").Value = "Tapu Fini"
$ws.Range("B26").Value = 0.06860000000000001
$ws.Range("A27").Value = "Blacephalon"
$ws.Range("B27").Value = 0.0668
$ws.Range("A28").Value = "Rotom-Wash"
$ws.Range("B28").Value = 0.06509999999999999
$ws.Range("A29").Value = "Dragonite"
$ws.Range("B29").Value = 0.0623
$ws.Range("A30").Value = "Mew"
$ws.Range("B30").Value = 0.0607
$ws.Range("A31").Value = "Hippowdon"
$ws.Range("B31").Value = 0.0591
$ws.Range("A32").Value = "Gastrodon"
$ws.Range("B32").Value = 0.0535
$ws.Range("A33").Value = "Volcarona"
$ws.Range("B33").Value = 0.0521
$ws.Range("A34").Value = "Bisharp"
$ws.Range("B34").Value = 0.0517
$ws.Range("A35").Value = "Skarmory"
$ws.Range("B35").Value = 0.049
$ws.Range("A36").Value = "Marowak-Alola"
$ws.Range("B36").Value = 0.04849999999999999
$ws.Range("A37").Value = "Magnezone"
$ws.Range("B37").Value = 0.0472
$ws.Range("A38").Value = "Drampa"
$ws.Range("B38").Value = 0.044
$ws.Range("A39").Value = "Buzzwole"
$ws.Range("B39").Value = 0.0426
$ws.Range("A40").Value = "Victini"
$ws.Range("B40").Value = 0.0394
$ws.Range("A41").Value = "Blissey"
$ws.Range("B41").Value = 0.0372
$ws.Range("A42").Value = "Pelipper"
$ws.Range("B42").Value = 0.0359
$ws.Range("A43").Value = "Stakataka"
$ws.Range("B43").Value = 0.0354
$ws.Range("A44").Value = "Scizor"
$ws.Range("B44").Value = 0.0352
$ws.Range("A45").Value = "Nidoking"
$ws.Range("B45").Value = 0.0338
$ws.Range("A46").Value = "Ninetales-Alola"
$ws.Range("B46").Value = 0.0335
$ws.Range("A47").Value = "Barraskewda"
$ws.Range("B47").Value = 0.0312
$ws.Range("A48").Value = "Volcanion"
$ws.Range("B48").Value = 0.0266
$ws.Range("A49").Value = "Slowking"
$ws.Range("B49").Value = 0.0253
$ws.Range("A50").Value = "Celesteela"
$ws.Range("B50").Value = 0.0232
$ws.Range("A51").Value = "Moltres-Galar"
$ws.Range("B51").Value = 0.0231
$ws.Range("A52").Value = "Aegislash"
$ws.Range("B52").Value = 0.023
$ws.Range("A53").Value = "Cloyster"
$ws.Range("B53").Value = 0.0225
$ws.Range("A54").Value = "Umbreon"
$ws.Range("B54").Value = 0.0222
$ws.Range("A55").Value = "Seismitoad"
$ws.Range("B55").Value = 0.0213
$ws.Range("A56").Value = "Dracozolt"
$ws.Range("B56").Value = 0.021
$ws.Range("A57").Value = "Swampert"
$ws.Range("B57").Value = 0.0204
$ws.Range("A58").Value = "Regieleki"
$ws.Range("B58").Value = 0.0202
$ws.Range("A59").Value = "Crawdaunt"
$ws.Range("B59").Value = 0.0201
$ws.Range("A60").Value = "Hawlucha"
$ws.Range("B60").Value = 0.0199
$ws.Range("A61").Value = "Hydreigon"
$ws.Range("B61").Value = 0.0164
$ws.Range("A62").Value = "Nihilego"
$ws.Range("B62").Value = 0.0163
$ws.Range("A63").Value = "Quagsire"
$ws.Range("B63").Value = 0.0153
$ws.Range("A64").Value = "Suicune"
$ws.Range("B64").Value = 0.0141
$ws.Range("A65").Value = "Arctozolt"
$ws.Range("B65").Value = 0.0124
$ws.Range("A66").Value = "Moltres"
$ws.Range("B66").Value = 0.0118
$ws.Range("A67").Value = "Porygon2"
$ws.Range("B67").Value = 0.0107
$ws.Range("A68").Value = "Zapdos-Galar"
$ws.Range("B68").Value = 0.0104
$ws.Range("A69").Value = "Reuniclus"
$ws.Range("B69").Value = 0.0103

# Rows 70-73 are no longer part of the dataset; removing them also shrinks the
# sheet dimension from A1:B73 down to A1:B69.
$ws.Rows("70:73").Delete()

# ColumnWidth is expressed in "characters"; Excel stores width internally with
# a fixed +5/6 padding offset, so 24.666... here serializes out as 25.5,
# matching column B's width exactly.
$ws.Range("A:A").ColumnWidth = 24.666666666666668
